# Updates cryptos list data (coin name, link, price, volume) for rows 2-51
# on Sheet1, matching a refreshed pull from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 2).Value = "Bitcoin"
$ws.Cells.Item(2, 3).Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Cells.Item(2, 4).Value = "25.962.00"
$ws.Cells.Item(2, 5).Value = "  +0.48%  "

# Row 3: Ethereum
$ws.Cells.Item(3, 2).Value = "Ethereum"
$ws.Cells.Item(3, 3).Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Cells.Item(3, 4).Value = "1.757.27"
$ws.Cells.Item(3, 5).Value = "  +0.49%  "

# Row 4: TetherUSD
$ws.Cells.Item(4, 2).Value = "TetherUSD"
$ws.Cells.Item(4, 3).Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9997"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.25%  "

# Row 5: BNB
$ws.Cells.Item(5, 2).Value = "BNB"
$ws.Cells.Item(5, 3).Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "237.09"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.22%  "

# Row 6: USDC
$ws.Cells.Item(6, 2).Value = "USDC"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.9993"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.21%  "

# Row 7: XRP
$ws.Cells.Item(7, 2).Value = "XRP"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.5218"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +3.85%  "

# Row 8: Cardano
$ws.Cells.Item(8, 2).Value = "Cardano"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2709"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +2.98%  "

# Row 9: Dogecoin
$ws.Cells.Item(9, 2).Value = "Dogecoin"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06213"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.95%  "

# Row 10: WrappedEther
$ws.Cells.Item(10, 2).Value = "WrappedEther"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(10, 4).Value = "1.757.75"
$ws.Cells.Item(10, 5).Value = "  +0.43%  "

# Row 11: TRON
$ws.Cells.Item(11, 2).Value = "TRON"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07026"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.21%  "

# Row 12: Solana
$ws.Cells.Item(12, 2).Value = "Solana"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "15.70"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +2.73%  "

# Row 13: Polygon
$ws.Cells.Item(13, 2).Value = "Polygon"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.6508"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +10.70%  "

# Row 14: Polkadot
$ws.Cells.Item(14, 2).Value = "Polkadot"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.495"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.33%  "

# Row 15: Litecoin
$ws.Cells.Item(15, 2).Value = "Litecoin"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "78.26"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.41%  "

# Row 16: BinanceUSD
$ws.Cells.Item(16, 2).Value = "BinanceUSD"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.9992"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.21%  "

# Row 17: Dai
$ws.Cells.Item(17, 2).Value = "Dai"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.9989"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.15%  "

# Row 18: WrappedBTC
$ws.Cells.Item(18, 2).Value = "WrappedBTC"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(18, 4).Value = "25.986.68"
$ws.Cells.Item(18, 5).Value = "  +0.20%  "

# Row 19: Avalanche
$ws.Cells.Item(19, 2).Value = "Avalanche"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "11.70"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.65%  "

# Row 20: ShibaInu
$ws.Cells.Item(20, 2).Value = "ShibaInu"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.000006712"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.51%  "

# Row 21: WrappedliquidstakedEther2.0
$ws.Cells.Item(21, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(21, 4).Value = "1.979.99"
$ws.Cells.Item(21, 5).Value = "  +0.33%  "

# Row 22: Uniswap
$ws.Cells.Item(22, 2).Value = "Uniswap"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.085"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.67%  "

# Row 23: Cosmos
$ws.Cells.Item(23, 2).Value = "Cosmos"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "8.411"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +4.04%  "

# Row 24: Chainlink
$ws.Cells.Item(24, 2).Value = "Chainlink"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "5.212"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +2.11%  "

# Row 25: Monero
$ws.Cells.Item(25, 2).Value = "Monero"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "137.11"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.66%  "

# Row 26: Toncoin
$ws.Cells.Item(26, 2).Value = "Toncoin"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "1.483"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -3.78%  "

# Row 27: EthereumClassic
$ws.Cells.Item(27, 2).Value = "EthereumClassic"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "15.19"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +2.30%  "

# Row 28: LidoDAOToken
$ws.Cells.Item(28, 2).Value = "LidoDAOToken"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.820"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -0.80%  "

# Row 29: BitcoinCash
$ws.Cells.Item(29, 2).Value = "BitcoinCash"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "103.27"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.58%  "

# Row 30: Stellar
$ws.Cells.Item(30, 2).Value = "Stellar"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.08402"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +3.47%  "

# Row 31: InternetComputer(DFINITY)
$ws.Cells.Item(31, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "3.700"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -1.58%  "

# Row 32: Filecoin
$ws.Cells.Item(32, 2).Value = "Filecoin"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.420"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -1.09%  "

# Row 33: Hedera
$ws.Cells.Item(33, 2).Value = "Hedera"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.04422"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -1.17%  "

# Row 34: HuobiToken
$ws.Cells.Item(34, 2).Value = "HuobiToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.646"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +0.14%  "

# Row 35: ARBITRUM
$ws.Cells.Item(35, 2).Value = "ARBITRUM"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.001"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +1.20%  "

# Row 36: ImmutableX
$ws.Cells.Item(36, 2).Value = "ImmutableX"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.6098"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +1.23%  "

# Row 37: MXToken
$ws.Cells.Item(37, 2).Value = "MXToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.734"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +2.29%  "

# Row 38: VeChain
$ws.Cells.Item(38, 2).Value = "VeChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.01570"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +2.29%  "

# Row 39: RenderToken
$ws.Cells.Item(39, 2).Value = "RenderToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.953"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.88%  "

# Row 40: PaxDollar
$ws.Cells.Item(40, 2).Value = "PaxDollar"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.9996"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.00%  "

# Row 41: Quant
$ws.Cells.Item(41, 2).Value = "Quant"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "102.80"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.90%  "

# Row 42: TheSandbox
$ws.Cells.Item(42, 2).Value = "TheSandbox"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.3895"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +2.39%  "

# Row 43: TrustWalletToken
$ws.Cells.Item(43, 2).Value = "TrustWalletToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.7510"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +2.92%  "

# Row 44: FraxShare
$ws.Cells.Item(44, 2).Value = "FraxShare"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "4.923"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -4.25%  "

# Row 45: Cronos
$ws.Cells.Item(45, 2).Value = "Cronos"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.05495"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +3.64%  "

# Row 46: Algorand
$ws.Cells.Item(46, 2).Value = "Algorand"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.1119"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.20%  "

# Row 47: Aptos
$ws.Cells.Item(47, 2).Value = "Aptos"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "6.109"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +2.93%  "

# Row 48: Elrond
$ws.Cells.Item(48, 2).Value = "Elrond"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "30.25"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +0.47%  "

# Row 49: Aave
$ws.Cells.Item(49, 2).Value = "Aave"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "52.83"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.92%  "

# Row 50: USDD
$ws.Cells.Item(50, 2).Value = "USDD"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.002"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.15%  "

# Row 51: Decentraland
$ws.Cells.Item(51, 2).Value = "Decentraland"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.3442"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.20%  "

